$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the selected cell / active cell in the sheet view
$ws.Range("G13").Select()

# G8 = G7 - 15
$ws.Range("G8").Formula = "=G7- 15"

# G10 = G9 - 15
$ws.Range("G10").Formula = "=G9- 15"

# G13 = G12 - 15
$ws.Range("G13").Formula = "=G12- 15"
